$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2023" column (L) to the right of the existing "2022" column (K),
# mirroring the formatting of column K for the header-border row (3),
# the year-label row (4) and the data row (5).

$ws.Range("K3").Copy() | Out-Null
$ws.Range("L3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("K4").Copy() | Out-Null
$ws.Range("L4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("K5").Copy() | Out-Null
$ws.Range("L5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# New values for the 2023 column
$ws.Range("L4").Value2 = 2023
$ws.Range("L5").Value2 = 0.11972285283622097

# Row 5 grows slightly taller to fit the wrapped text with the extra column
$ws.Rows(5).RowHeight = 40.5

# Reset selection back to the top-left cell
$ws.Range("A1").Select() | Out-Null
